# Update gh-pages to output generated at 456a3b4
# Applies refreshed "want to go" counts (column F) and one venue name
# update (column C) on the "展览" (sheet 1) and "全部类型" (sheet 4)
# worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- Sheet "展览" ----
$ws1.Range("F3").Value = 30
$ws1.Range("F4").Value = 21
$ws1.Range("F5").Value = 16065
$ws1.Range("F7").Value = 12
$ws1.Range("F8").Value = 723
$ws1.Range("F9").Value = 15532
$ws1.Range("F10").Value = 64
$ws1.Range("F11").Value = 9152
$ws1.Range("F12").Value = 410
$ws1.Range("F13").Value = 10
$ws1.Range("F14").Value = 1024
$ws1.Range("F15").Value = 113
$ws1.Range("F19").Value = 21
$ws1.Range("F20").Value = 73
$ws1.Range("F21").Value = 579
$ws1.Range("F22").Value = 26
$ws1.Range("F23").Value = 13
$ws1.Range("F24").Value = 69
$ws1.Range("F26").Value = 12
$ws1.Range("F27").Value = 19
$ws1.Range("F28").Value = 29
$ws1.Range("F29").Value = 509
$ws1.Range("F33").Value = 72
$ws1.Range("F34").Value = 60
$ws1.Range("F36").Value = 339
$ws1.Range("F37").Value = 467
$ws1.Range("F39").Value = 5632
$ws1.Range("C40").Value = "苏州·星部落&青铜树动漫嘉年华"
$ws1.Range("F40").Value = 5238

# ---- Sheet "全部类型" ----
$ws4.Range("F3").Value = 30
$ws4.Range("F4").Value = 21
$ws4.Range("F5").Value = 16066
$ws4.Range("F7").Value = 12
$ws4.Range("F8").Value = 723
$ws4.Range("F9").Value = 15532
$ws4.Range("F10").Value = 64
$ws4.Range("F11").Value = 9152
$ws4.Range("F12").Value = 410
$ws4.Range("F13").Value = 10
$ws4.Range("F14").Value = 1024
$ws4.Range("F15").Value = 113
$ws4.Range("F19").Value = 21
$ws4.Range("F20").Value = 73
$ws4.Range("F21").Value = 579
$ws4.Range("F22").Value = 26
$ws4.Range("F23").Value = 13
$ws4.Range("F24").Value = 69
$ws4.Range("F26").Value = 12
$ws4.Range("F27").Value = 19
$ws4.Range("F28").Value = 29
$ws4.Range("F29").Value = 509
$ws4.Range("F35").Value = 72
$ws4.Range("F36").Value = 60
$ws4.Range("F38").Value = 339
$ws4.Range("F39").Value = 467
$ws4.Range("F41").Value = 5632
$ws4.Range("C43").Value = "苏州·星部落&青铜树动漫嘉年华"
$ws4.Range("F43").Value = 5238
